# Applies the row-content swaps described by the diff.
# Pattern: several pairs of rows had their entire (non-row-number) content
# swapped between them (species records re-ordered), while row numbers
# (the "r" attribute / sheet row position) themselves stayed put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 6 <-> 8 -----------------------------------------------------
$ws.Range("A6").Value = 131066788
$ws.Range("B6").Value = 83215
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 308
$ws.Range("F6").Value = "Brunpudrad nållav"
$ws.Range("G6").Value = "Chaenotheca gracillima"
$ws.Range("H6").Value = "(Vain.) Tibell"
$ws.Range("Q6").Value = 425211
$ws.Range("R6").Value = 6712276

$ws.Range("A8").Value = 131066787
$ws.Range("B8").Value = 80383
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 6463
$ws.Range("F8").Value = "Bårdlav"
$ws.Range("G8").Value = "Nephroma parile"
$ws.Range("H8").Value = "(Ach.) Ach."
$ws.Range("Q8").Value = 425069
$ws.Range("R8").Value = 6712290

# --- Rows 9 <-> 10 ------------------------------------------------------
# (row 9 also loses its K/L/M/N/AC content, which moves onto row 10)
$ws.Range("A9").Value = 131066783
$ws.Range("B9").Value = 83089
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 1312
$ws.Range("F9").Value = "Gammelgransskål"
$ws.Range("G9").Value = "Pseudographis pinicola"
$ws.Range("H9").Value = "(Nyl.) Rehm"
$ws.Range("Q9").Value = 425170
$ws.Range("R9").Value = 6712292
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("AC9").ClearContents()

$ws.Range("A10").Value = 131066770
$ws.Range("B10").Value = 57884
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 100109
$ws.Range("F10").Value = "Tretåig hackspett"
$ws.Range("G10").Value = "Picoides tridactylus"
$ws.Range("H10").Value = "(Linnaeus, 1758)"
$ws.Range("Q10").Value = 425323
$ws.Range("R10").Value = 6712206
$ws.Range("K10").Value = ""
$ws.Range("L10").Value = ""
$ws.Range("M10").Value = "äldre spår"
$ws.Range("N10").Value = ""
$ws.Range("AC10").Value = "Ringhack på gran"

# --- Rows 15 <-> 17 -------------------------------------------------
$ws.Range("A15").Value = 131066761
$ws.Range("B15").Value = 91771
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 5447
$ws.Range("F15").Value = "Vedticka"
$ws.Range("G15").Value = "Fuscoporia viticola"
$ws.Range("H15").Value = "(Schwein.) Murrill"
$ws.Range("Q15").Value = 425072
$ws.Range("R15").Value = 6712273

$ws.Range("A17").Value = 131066768
$ws.Range("B17").Value = 91808
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 1202
$ws.Range("F17").Value = "Ullticka"
$ws.Range("G17").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H17").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q17").Value = 425256
$ws.Range("R17").Value = 6712203

# --- Rows 21 <-> 22 -------------------------------------------------
$ws.Range("A21").Value = 131066766
$ws.Range("B21").Value = 92179
$ws.Range("D21").Value = "VU"
$ws.Range("E21").Value = 2062
$ws.Range("F21").Value = "Ulltickeporing"
$ws.Range("G21").Value = "Skeletocutis brevispora"
$ws.Range("H21").Value = "Niemelä"
$ws.Range("Q21").Value = 425069
$ws.Range("R21").Value = 6712285

$ws.Range("A22").Value = 131066778
$ws.Range("B22").Value = 81228
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 1049
$ws.Range("F22").Value = "Kortskaftad ärgspik"
$ws.Range("G22").Value = "Microcalicium ahlneri"
$ws.Range("H22").Value = "Tibell"
$ws.Range("Q22").Value = 425336
$ws.Range("R22").Value = 6712202

# --- Rows 25 <-> 26 (only the Id and Nord columns differed) --------
$ws.Range("A25").Value = 131066762
$ws.Range("R25").Value = 6712254

$ws.Range("A26").Value = 131066765
$ws.Range("R26").Value = 6712253

# --- Rows 32 <-> 33 -------------------------------------------------
$ws.Range("A32").Value = 131066790
$ws.Range("B32").Value = 83215
$ws.Range("D32").Value = "NT"
$ws.Range("E32").Value = 308
$ws.Range("F32").Value = "Brunpudrad nållav"
$ws.Range("G32").Value = "Chaenotheca gracillima"
$ws.Range("H32").Value = "(Vain.) Tibell"
$ws.Range("Q32").Value = 425164
$ws.Range("R32").Value = 6712278

$ws.Range("A33").Value = 131066767
$ws.Range("B33").Value = 91808
$ws.Range("D33").Value = "NT"
$ws.Range("E33").Value = 1202
$ws.Range("F33").Value = "Ullticka"
$ws.Range("G33").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H33").Value = "(Ach.) Ach."
$ws.Range("Q33").Value = 425259
$ws.Range("R33").Value = 6712201
